# Fruta / hortaliza, semanal
# Insert a new weekly record at row 191 (pushing the existing rows 191-215
# down to 192-216), then populate the new row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 191; this shifts rows 191-215
# down to 192-216 and extends the sheet dimension to A1:R216 automatically.
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row 191 with the new record's data.
$ws.Range("A191").Value = 10
$ws.Range("B191").Value = "Vega Modelo de Temuco"
$ws.Range("C191").Value = "La Araucanía"
$ws.Range("D191").Value = 44449
$ws.Range("E191").Value = 9
$ws.Range("F191").Value = 100112008
$ws.Range("G191").Value = "Coliflor"
$ws.Range("H191").Value = "Sin especificar"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 2350
$ws.Range("K191").Value = 800
$ws.Range("L191").Value = 900
$ws.Range("M191").Value = 853
$ws.Range("N191").Value = "$/unidad"
$ws.Range("O191").Value = "Región Metropolitana"
$ws.Range("P191").Value = 853
$ws.Range("Q191").Value = 1
$ws.Range("R191").Value = "Hortaliza"
